# "Colocando header nos graficos" - add a header label to column A (row 1)
# on each chart-data sheet, un-bold/un-border the existing row labels
# (A2:A12), and fix accented Portuguese spelling of several labels.
# Also trims the "Teto" row from the emissions sheet and refreshes the
# cost sheet header/values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheets 1-4 share the exact same row layout:
#   row 1:  (blank) | 2015 | 2030 | 2040 | 2050
#   row 2:  Hidro
#   row 3:  Gas Natural      -> Gás Natural
#   row 4:  Carvao           -> Carvão
#   row 5:  Nuclear
#   row 6:  Oleos Comb       -> Óleos Comb
#   row 7:  Biomassa
#   row 8:  Eolica           -> Eólica
#   row 9:  Solar
#   row 10: Outros
#   row 11: Pot Compl        -> Pot. Compl.
#   row 12: GD
# ---------------------------------------------------------------------

$labels = @{
    2  = "Hidro"
    3  = "Gás Natural"
    4  = "Carvão"
    5  = "Nuclear"
    6  = "Óleos Comb"
    7  = "Biomassa"
    8  = "Eólica"
    9  = "Solar"
    10 = "Outros"
    11 = "Pot. Compl."
    12 = "GD"
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # A1 header: copy style from B1 (bold/border/center) then set text
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # A2:A12 - drop the bold/border style and fix accented text
    foreach ($r in 2..12) {
        $cell = $ws.Cells.Item($r, 1)
        $cell.ClearFormats()
        $cell.Value = $labels[$r]
    }
}

# ---------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
#   row 1:  (blank) | 2015 | 2030 | 2040 | 2050
#   row 2:  P Medio   -> P.Médio
#   row 3:  P Critico -> P.Crítico
#   row 4:  Teto      -> removed entirely
# ---------------------------------------------------------------------

$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws5.Range("A1").Value = "Período"

$ws5.Cells.Item(2, 1).ClearFormats()
$ws5.Cells.Item(2, 1).Value = "P.Médio"

$ws5.Cells.Item(3, 1).ClearFormats()
$ws5.Cells.Item(3, 1).Value = "P.Crítico"

$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)"
#   row 1:  (blank) | Custo               -> (blank) | 2015
#   row 2:  Expansao Centralizada | 739   -> Expansão Centralizada | 586
#   row 3:  Expansao por GD | 65          -> Expansão por GD | 99
# ---------------------------------------------------------------------

$ws6 = $wb.Worksheets.Item(6)

# A1 header: copy style from B1 before B1's own text changes
$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws6.Range("A1").Value = "Tipo Expansão"

# B1 must stay a text value "2015" (not get auto-converted to a number)
$ws6.Range("B1").NumberFormat = "@"
$ws6.Range("B1").Value = "2015"
# Restore B1's original look (bold/border/center) using A1, which now
# carries the exact same style B1 originally had.
$ws6.Range("A1").Copy()
$ws6.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws6.Cells.Item(2, 1).ClearFormats()
$ws6.Cells.Item(2, 1).Value = "Expansão Centralizada"
$ws6.Cells.Item(2, 2).Value = 586

$ws6.Cells.Item(3, 1).ClearFormats()
$ws6.Cells.Item(3, 1).Value = "Expansão por GD"
$ws6.Cells.Item(3, 2).Value = 99
